{"js": "// Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific bullet /\n// impact paragraphs, per the commit's \"quantitative metrics highlighting\"\n// change. Each target paragraph's plain run is located by its exact\n// original text, then the numeric substrings inside it are located with\n// a paragraph-scoped search and given bold + color formatting \u2014 Word\n// automatically splits the run(s) around the match, preserving\n// whitespace via xml:space=\"preserve\" where needed.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Ordered list of [originalParagraphText, [metricsToHighlight...]]\nconst EDITS = [\n  [\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    [\"23%\", \"64%\"],\n  ],\n  [\n    \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n    [\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\"],\n  ],\n  [\n    \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    [\"1,200\"],\n  ],\n  [\n    \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    [\"$400M\", \"$1B\"],\n  ],\n  [\n    \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    [\"73.5%\", \"$4.7M\"],\n  ],\n  [\n    \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    [\"87%\", \"71%\"],\n  ],\n];\n\n// Load all paragraphs + their text once up front.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\n// Map each target paragraph text to its Paragraph object. Track how many\n// times each target text has already been consumed, so duplicate texts\n// (e.g. the two \"Achieved 87% ... 71%\" variants) resolve to distinct\n// paragraphs in document order, matching the diff's per-occurrence edits.\nconst consumedCount = new Map();\n\nfor (const [paragraphText, metrics] of EDITS) {\n  const already = consumedCount.get(paragraphText) || 0;\n  let seen = 0;\n  let targetParagraph = null;\n\n  for (const p of paragraphs.items) {\n    if (p.text === paragraphText) {\n      if (seen === already) {\n        targetParagraph = p;\n        break;\n      }\n      seen++;\n    }\n  }\n\n  if (!targetParagraph) {\n    throw new Error(\"Could not locate paragraph: \" + paragraphText);\n  }\n  consumedCount.set(paragraphText, already + 1);\n\n  // Highlight each metric substring within this paragraph, in order.\n  for (const metric of metrics) {\n    const scopedRange = targetParagraph.getRange();\n    const searchResults = scopedRange.search(metric, { matchCase: true });\n    searchResults.load(\"items\");\n    await context.sync();\n\n    if (searchResults.items.length === 0) {\n      throw new Error(\n        \"Could not locate metric '\" + metric + \"' in paragraph: \" + paragraphText\n      );\n    }\n\n    const matchRange = searchResults.items[0];\n    matchRange.font.bold = true;\n    matchRange.font.color = HIGHLIGHT_COLOR;\n    await context.sync();\n  }\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific bullet /\n# impact paragraphs, per the commit's \"quantitative metrics highlighting\"\n# change. Each target paragraph is located by its exact original text\n# (trimming the trailing paragraph-mark CR that Word COM includes), then\n# each numeric substring inside it is located with Find.Execute scoped to\n# that paragraph's Range and given bold + color formatting. Word\n# automatically splits the run(s) around the match, preserving whitespace\n# via xml:space=\"preserve\" where needed.\n\n$d = $word.ActiveDocument\n\n# #2C3E50 as a Word \"wdColor\" BGR-packed long: B*65536 + G*256 + R\n$HighlightColor = 0x50 * 65536 + 0x3E * 256 + 0x2C\n\n# Ordered list of (original paragraph text, metrics to highlight within it)\n$Edits = @(\n    @{\n        Text    = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Text    = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\"\n        Metrics = @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\")\n    },\n    @{\n        Text    = \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Metrics = @(\"1,200\")\n    },\n    @{\n        Text    = \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Metrics = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Text    = \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Metrics = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text    = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Metrics = @(\"87%\", \"71%\")\n    }\n)\n\n# Track how many times each target text has already been consumed, so\n# duplicate texts (e.g. the two \"Achieved 87% ... 71%\" variants) resolve\n# to distinct paragraphs in document order, matching the diff's\n# per-occurrence edits.\n$consumedCount = @{}\n\n$paragraphCount = $d.Paragraphs.Count\n\nforeach ($edit in $Edits) {\n    $wantText = $edit.Text\n    $already = 0\n    if ($consumedCount.ContainsKey($wantText)) {\n        $already = $consumedCount[$wantText]\n    }\n\n    $seen = 0\n    $targetParagraph = $null\n\n    for ($i = 1; $i -le $paragraphCount; $i++) {\n        $para = $d.Paragraphs.Item($i)\n        $paraText = $para.Range.Text.TrimEnd(\"`r\")\n        if ($paraText -eq $wantText) {\n            if ($seen -eq $already) {\n                $targetParagraph = $para\n                break\n            }\n            $seen = $seen + 1\n        }\n    }\n\n    if ($null -eq $targetParagraph) {\n        throw \"Could not locate paragraph: $wantText\"\n    }\n    $consumedCount[$wantText] = $already + 1\n\n    foreach ($metric in $edit.Metrics) {\n        $searchRange = $targetParagraph.Range\n        $found = $searchRange.Find.Execute($metric, $true)\n        if (-not $found) {\n            throw \"Could not locate metric '$metric' in paragraph: $wantText\"\n        }\n        $searchRange.Font.Bold = 1\n        $searchRange.Font.Color = $HighlightColor\n    }\n}\n"}
